# Applies the recorded field-level edits to the "Artfynd" worksheet.
# Rows 3-5 have most of their observation data cyclically rotated,
# rows 22-23 swap their species records (with distinct new B values),
# and a handful of other rows only get their "B" (sort order) value bumped by 1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 91809

# Row 3
$ws.Range("A3").Value = 130825852
$ws.Range("P3").Value = 'Flinktorpet, Kälom, Offerdal, Jmt'
$ws.Range("Q3").Value = 460952
$ws.Range("R3").Value = 7039723
$ws.Range("S3").Value = 15
$ws.Range("Z3").Value = '10:42'
$ws.Range("AB3").Value = '10:42'
$ws.Range("AC3").Value = 'Barkfläkta grövre och klenare granar.'

# Row 4
$ws.Range("A4").Value = 130825823
$ws.Range("B4").Value = 57881
$ws.Range("E4").Value = 100049
$ws.Range("F4").Value = 'Spillkråka'
$ws.Range("G4").Value = 'Dryocopus martius'
$ws.Range("M4").Value = 'äldre spår'
$ws.Range("Q4").Value = 460947
$ws.Range("R4").Value = 7039711
$ws.Range("S4").Value = 10
$ws.Range("Z4").Value = '10:38'
$ws.Range("AB4").Value = '10:38'
$ws.Range("AC4").Value = 'Födosökshål på äldre döende gran.'

# Row 5
$ws.Range("A5").Value = 130826784
$ws.Range("B5").Value = 57884
$ws.Range("E5").Value = 100109
$ws.Range("F5").Value = 'Tretåig hackspett'
$ws.Range("G5").Value = 'Picoides tridactylus'
$ws.Range("M5").Value = 'färska spår'
$ws.Range("P5").Value = 'Brännan, Kälom, Offerdal, Jmt'
$ws.Range("Q5").Value = 461233
$ws.Range("R5").Value = 7039438
$ws.Range("Z5").Value = '11:37'
$ws.Range("AB5").Value = '11:37'
$ws.Range("AC5").Value = 'Födosök barkfläk'

# Row 6
$ws.Range("B6").Value = 92268

# Row 8
$ws.Range("B8").Value = 79244

# Row 11
$ws.Range("B11").Value = 79244

# Row 12
$ws.Range("B12").Value = 91809

# Row 15
$ws.Range("B15").Value = 91809

# Row 16
$ws.Range("B16").Value = 91809

# Row 17
$ws.Range("B17").Value = 89194

# Row 18
$ws.Range("B18").Value = 91809

# Row 19
$ws.Range("B19").Value = 79244

# Row 22
$ws.Range("A22").Value = 130826438
$ws.Range("B22").Value = 79244
$ws.Range("D22").Value = 'NT'
$ws.Range("E22").Value = 6425
$ws.Range("F22").Value = 'Garnlav'
$ws.Range("G22").Value = 'Alectoria sarmentosa'
$ws.Range("H22").Value = '(Ach.) Ach.'
$ws.Range("P22").Value = 'Brännan, Brännan, Jmt'
$ws.Range("Q22").Value = 461220
$ws.Range("R22").Value = 7039590
$ws.Range("S22").Value = 25
$ws.Range("Z22").Value = '11:16'
$ws.Range("AB22").Value = '11:16'
$ws.Range("AC22").Value = 'Rikligt i området'

# Row 23
$ws.Range("A23").Value = 130826355
$ws.Range("B23").Value = 92536
$ws.Range("D23").Value = 'VU'
$ws.Range("E23").Value = 67
$ws.Range("F23").Value = 'Sprickporing'
$ws.Range("G23").Value = 'Diplomitoporus crustulinus'
$ws.Range("H23").Value = '(Bres.) Domański'
$ws.Range("P23").Value = 'Flinktorpet, Flinktorpet, Jmt'
$ws.Range("Q23").Value = 461117
$ws.Range("R23").Value = 7039629
$ws.Range("S23").Value = 10
$ws.Range("Z23").Value = '11:10'
$ws.Range("AB23").Value = '11:10'
$ws.Range("AC23").Value = 'På undersidan av lutande död gran.'
